$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates -------------------------------------------------

$ws.Range("A1").Value = "Item Name (Michael Edgar)"

$ws.Range("B2").Value = "The role of the game manager is to keep track of the current score, high score and lives count"

$ws.Range("B4").Value = "Communicate with the Pop up Score to update the current score, the world when the player"
$ws.Range("B5").Value = "dies to update the current lives and the update the high score when necessary"

$ws.Range("B9").Value = "Screen display of score, high score and lives"
$ws.Range("D9").Value = "Communication with?"
$ws.Range("E9").ClearContents()

$ws.Range("D10").Value = "Display/UI"

$ws.Range("B13").Value = "Receive score update from the pop up score when the player interacts with an item and receive"
$ws.Range("D13").Value = "Communication with?"
$ws.Range("E13").ClearContents()

$ws.Range("B14").Value = "update when player dies"
$ws.Range("D14").Value = "Pop up score and player"

# --- Column widths / visibility -------------------------------------------

$ws.Columns.Item(1).ColumnWidth = 47.83333333333333
$ws.Columns.Item(2).ColumnWidth = 84.66666666666667
$ws.Columns.Item(3).ColumnWidth = 87.66666666666667
$ws.Columns.Item(4).ColumnWidth = 38.666666666666664

$ws.Columns.Item(5).ColumnWidth = 9.833333333333332
$ws.Columns.Item(5).Hidden = $true

$ws.Columns.Item(6).ColumnWidth = 40.166666666666664
$ws.Columns.Item(6).Hidden = $true
$ws.Columns.Item(7).ColumnWidth = 40.166666666666664
$ws.Columns.Item(7).Hidden = $true

for ($col = 8; $col -le 26; $col++) {
    $ws.Columns.Item($col).ColumnWidth = 8.333333333333332
    $ws.Columns.Item($col).Hidden = $true
}

# --- Selection (matches last-active cell recorded in the saved file) ------

$ws.Range("D14").Select()
